# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览": 草莓动漫节 1218 -> 1221, 第一届ANE·DACG动漫嘉年华 610 -> 611
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1221
$wsExhibit.Range("F5").Value = 611

# Sheet "全部类型": same two events appear again with identical updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1221
$wsAll.Range("F6").Value = 611
